$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = Get-Date -Year 2020 -Month 4 -Day 29 -Hour 0 -Minute 0 -Second 0

# Row 21: rename "Teste" to "Leo" and update the "DATA DE SAIDA" date
$ws.Range("A21").Value = "Leo"
$ws.Range("C21").Value = $newDate

# Row 22: update the "DATA DE SAIDA" date for "Rodrigo " and match the date
# number format used elsewhere in the column (same style as C21) by copying
# C21's formatting onto C22
$ws.Range("C22").Value = $newDate
$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active selection to reflect the cell the author was last on
$ws.Range("C22").Select()
